$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("09/11/2023")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("6. Devise a scheme in computing a polynomial ‘C’ using arrays, where ‘C’ is computed by:")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(a) Adding two polynomials A and B")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(b) Subtracting polynomial B from polynomial A")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(c) Multiplying two polynomials A and B")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(d) Differentiating polynomial A")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("7. Devise a scheme to represent a sparse matrix X and transpose this representation of X in lexicographic order.")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("23/11/2023")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("8. Implement Assignment (6) using dynamic data structure as follows:")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("Devise schemes of dynamic data structures to compute a polynomial ‘C’ where ‘C’ is computed by:")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(a) Adding two polynomials A and B,")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(b) Subtracting polynomial B from polynomial A,")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(c) Multiplying two polynomials A and B,")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("(d) Differentiating polynomial A.")

